$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural edit: insert a new column before the old "E" (Description) column.
# This shifts the old column E (with its "Description"-header cell style) to F,
# and creates a brand-new column E that inherits the row/column formatting.
$ws.Columns("E:E").Insert()

# Remove the old "Project 2" / "Project 3" sample rows (rows 3 and 4); the remaining
# "Project 4" sample row (old row 5) shifts up to become row 3, and the two blank
# trailer rows (old rows 6-7) shift up to rows 4-5.
$ws.Rows("3:4").Delete()

# --- Header row ---
$ws.Range("E1").Value = "SRA"

# --- Row 2: rename "Project 1" -> "Project PE" and give it its own description ---
$ws.Range("A2").Value = "Project PE"
$ws.Range("F2").Value = "test batch with PE project"

# --- Row 3 (previously "Project 4"): rename -> "Project SE", new description ---
$ws.Range("A3").Value = "Project SE"
$ws.Range("F3").Value = "test batch with SE input"

# --- Row 4 (previously blank): new SRA-accession sample row ---
$ws.Range("A4").Value = "SRR11241255"
$ws.Range("E4").Value = "SRR11241255"
$ws.Range("F4").Value = "test batch with SRA accession"

# --- Column widths ---
$ws.Columns("A:A").ColumnWidth = 12.33203125
$ws.Columns("E:E").ColumnWidth = 22.1640625
$ws.Columns("F:F").ColumnWidth = 24.83203125

# --- New cell comments describing the new D1/E1 headers ---
$d1Comment = $ws.Range("D1").AddComment("Chienchi Lo:" + [char]10 + "Single End Reads")
$e1Comment = $ws.Range("E1").AddComment("Chienchi Lo:" + [char]10 + "Input SRA accession")

# --- Selection matches the saved view in the edited workbook ---
$ws.Range("D15").Select()
